$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Agrn"
$ws.Range("C2").Value = "Musk"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 15.59657533333333
$ws.Range("H2").Value = 46.789726
$ws.Range("I2").Value = 0.4757744772251148
$ws.Range("J2").Value = 0.475774477225115
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3212813333333333
$ws.Range("N2").Value = 0.9638439999999999
$ws.Range("O2").Value = 0.02198075961745464
$ws.Range("P2").Value = 0.02198075961745463
$ws.Range("Q2").Value = 5.010888518527111
$ws.Range("R2").Value = 45.097996666744
$ws.Range("S2").Value = 0.01045788441600539
$ws.Range("T2").Value = 0.0104578844160054

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Agrn"
$ws.Range("C3").Value = "Musk"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.59657533333333
$ws.Range("H3").Value = 46.789726
$ws.Range("I3").Value = 0.4757744772251148
$ws.Range("J3").Value = 0.475774477225115
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.908863
$ws.Range("N3").Value = 26.726589
$ws.Range("O3").Value = 0.6095081031821615
$ws.Range("P3").Value = 0.6095081031821615
$ws.Range("Q3").Value = 138.947752913846
$ws.Range("R3").Value = 1250.529776224614
$ws.Range("S3").Value = 0.2899883991559642
$ws.Range("T3").Value = 0.2899883991559643

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Agrn"
$ws.Range("C4").Value = "Musk"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 15.59657533333333
$ws.Range("H4").Value = 46.789726
$ws.Range("I4").Value = 0.4757744772251148
$ws.Range("J4").Value = 0.475774477225115
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.386335666666667
$ws.Range("N4").Value = 16.159007
$ws.Range("O4").Value = 0.3685111372003839
$ws.Range("P4").Value = 0.3685111372003838
$ws.Range("Q4").Value = 84.00838999578689
$ws.Range("R4").Value = 756.075509962082
$ws.Range("S4").Value = 0.1753281936531452
$ws.Range("T4").Value = 0.1753281936531452

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Agrn"
$ws.Range("C5").Value = "Musk"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.399531333333333
$ws.Range("H5").Value = 19.198594
$ws.Range("I5").Value = 0.1952180917624358
$ws.Range("J5").Value = 0.1952180917624358
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.3212813333333333
$ws.Range("N5").Value = 0.9638439999999999
$ws.Range("O5").Value = 0.02198075961745464
$ws.Range("P5").Value = 0.02198075961745463
$ws.Range("Q5").Value = 2.056049959481777
$ws.Range("R5").Value = 18.504449635336
$ws.Range("S5").Value = 0.004291041948008301
$ws.Range("T5").Value = 0.004291041948008301

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Agrn"
$ws.Range("C6").Value = "Musk"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.399531333333333
$ws.Range("H6").Value = 19.198594
$ws.Range("I6").Value = 0.1952180917624358
$ws.Range("J6").Value = 0.1952180917624358
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.908863
$ws.Range("N6").Value = 26.726589
$ws.Range("O6").Value = 0.6095081031821615
$ws.Range("P6").Value = 0.6095081031821615
$ws.Range("Q6").Value = 57.012547912874
$ws.Range("R6").Value = 513.112931215866
$ws.Range("S6").Value = 0.1189870088169634
$ws.Range("T6").Value = 0.1189870088169634

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Agrn"
$ws.Range("C7").Value = "Musk"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.399531333333333
$ws.Range("H7").Value = 19.198594
$ws.Range("I7").Value = 0.1952180917624358
$ws.Range("J7").Value = 0.1952180917624358
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.386335666666667
$ws.Range("N7").Value = 16.159007
$ws.Range("O7").Value = 0.3685111372003839
$ws.Range("P7").Value = 0.3685111372003838
$ws.Range("Q7").Value = 34.47002387068422
$ws.Range("R7").Value = 310.230214836158
$ws.Range("S7").Value = 0.0719400409974641
$ws.Range("T7").Value = 0.0719400409974641

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Agrn"
$ws.Range("C8").Value = "Musk"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.78533933333333
$ws.Range("H8").Value = 32.356018
$ws.Range("I8").Value = 0.3290074310124493
$ws.Range("J8").Value = 0.3290074310124493
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.3212813333333333
$ws.Range("N8").Value = 0.9638439999999999
$ws.Range("O8").Value = 0.02198075961745464
$ws.Range("P8").Value = 0.02198075961745463
$ws.Range("Q8").Value = 3.465128201465777
$ws.Range("R8").Value = 31.186153813192
$ws.Range("S8").Value = 0.007231833253440937
$ws.Range("T8").Value = 0.007231833253440937

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Agrn"
$ws.Range("C9").Value = "Musk"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.78533933333333
$ws.Range("H9").Value = 32.356018
$ws.Range("I9").Value = 0.3290074310124493
$ws.Range("J9").Value = 0.3290074310124493
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 8.908863
$ws.Range("N9").Value = 26.726589
$ws.Range("O9").Value = 0.6095081031821615
$ws.Range("P9").Value = 0.6095081031821615
$ws.Range("Q9").Value = 96.085110529178
$ws.Range("R9").Value = 864.765994762602
$ws.Range("S9").Value = 0.2005326952092338
$ws.Range("T9").Value = 0.2005326952092339

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Agrn"
$ws.Range("C10").Value = "Musk"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 10.78533933333333
$ws.Range("H10").Value = 32.356018
$ws.Range("I10").Value = 0.3290074310124493
$ws.Range("J10").Value = 0.3290074310124493
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.386335666666667
$ws.Range("N10").Value = 16.159007
$ws.Range("O10").Value = 0.3685111372003839
$ws.Range("P10").Value = 0.3685111372003838
$ws.Range("Q10").Value = 58.09345792823622
$ws.Range("R10").Value = 522.8411213541259
$ws.Range("S10").Value = 0.1212429025497745
$ws.Range("T10").Value = 0.1212429025497745
